# Applies the diff described in the commit:
#   "progress: find bks, comparison table and statistical analysis working"
#
# Sheet "Resumen": B2 Zona label + C2 Maximo value updated.
# Sheet "Solucion": the Pedido/Salida assignment table is re-shuffled.
# Sheet "Metricas": Z1/Z2 Tiempo values updated.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Resumen" ----
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = "Z1"
$wsResumen.Range("C2").Value = 514.3404476159586

# ---- Sheet "Solucion" ----
$wsSolucion = $wb.Worksheets.Item("Solucion")

$pedidos = @(
    "Pedido_28", "Pedido_6", "Pedido_31", "Pedido_17", "Pedido_1",
    "Pedido_24", "Pedido_22", "Pedido_5", "Pedido_36", "Pedido_38",
    "Pedido_15", "Pedido_39", "Pedido_26", "Pedido_20", "Pedido_4",
    "Pedido_16", "Pedido_12", "Pedido_18", "Pedido_35", "Pedido_9",
    "Pedido_11", "Pedido_14", "Pedido_7", "Pedido_3", "Pedido_2",
    "Pedido_34", "Pedido_25", "Pedido_27", "Pedido_33", "Pedido_21",
    "Pedido_19", "Pedido_32", "Pedido_40", "Pedido_30", "Pedido_23",
    "Pedido_10", "Pedido_13", "Pedido_8", "Pedido_37", "Pedido_29"
)

$salidas = @(
    "S001", "S021", "S011", "S031", "S002",
    "S022", "S012", "S032", "S003", "S023",
    "S013", "S033", "S024", "S004", "S034",
    "S014", "S025", "S005", "S035", "S015",
    "S026", "S006", "S036", "S016", "S027",
    "S007", "S017", "S037", "S008", "S028",
    "S038", "S018", "S009", "S029", "S019",
    "S039", "S010", "S030", "S040", "S020"
)

for ($i = 0; $i -lt 40; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $pedidos[$i]
    $wsSolucion.Cells.Item($row, 2).Value = $salidas[$i]
}

# ---- Sheet "Metricas" ----
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 514.3404476159586
$wsMetricas.Range("B3").Value = 512.6677316466646
